$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column C (current RF column), shifting RF, XGBoost, KM, NN
# one column to the right to make room for the new "PR" (polynomial regression) column.
$ws.Range("C:C").Insert()

# Populate the new column C header / search type / parameter combinations.
$ws.Range("C1").Value = "PR"
$ws.Range("B1").Copy()
$ws.Range("C1").PasteSpecial(-4122)

$ws.Range("C2").Value = "GridSearchCV"
$ws.Range("C3").Value = 4

# Update the final computed timing results across the whole table (B4:G6).
$ws.Range("B4").Value = "54s"
$ws.Range("C4").Value = "32m, 19s"
$ws.Range("D4").Value = "15m, 45s"
$ws.Range("E4").Value = "13m, 12s"
$ws.Range("F4").Value = "47m, 36s"
$ws.Range("G4").Value = "1h, 26m"

$ws.Range("B5").Value = "1s"
$ws.Range("C5").Value = "10m, 41s"
$ws.Range("D5").Value = "1h, 1m"
$ws.Range("E5").Value = "49m, 3s"
$ws.Range("F5").Value = "2h, 15m"
$ws.Range("G5").Value = "5h, 41m"

$ws.Range("B6").Value = "3m, 2s"
$ws.Range("C6").Value = "4s"
$ws.Range("D6").Value = "15s"
$ws.Range("E6").Value = "1m, 52s"
$ws.Range("F6").Value = "50m, 7s"
$ws.Range("G6").Value = "6s"
